$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.349.71'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.936.99'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.16'
$ws.Range('E5').Value = '  +2.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7243'
$ws.Range('E6').Value = '  +3.69%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3313'
$ws.Range('E8').Value = '  +2.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '28.06'
$ws.Range('E9').Value = '  +6.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07273'
$ws.Range('E10').Value = '  +6.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8096'
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08106'
$ws.Range('E12').Value = '  +2.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.936.88'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.487'
$ws.Range('E14').Value = '  +2.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.95'
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.13'
$ws.Range('E16').Value = '  +4.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.342.34'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008258'
$ws.Range('E18').Value = '  +5.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '253.24'
$ws.Range('E19').Value = '  -2.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.844'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.191.28'
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.964'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.772'
$ws.Range('E25').Value = '  +1.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.08'
$ws.Range('E26').Value = '  +4.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.352'
$ws.Range('E27').Value = '  +5.92%  '
$ws.Range('E28').Value = '  +3.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1300'
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.544'
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.447'
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.215'
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05257'
$ws.Range('E34').Value = '  +4.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.268'
$ws.Range('E35').Value = '  +6.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7520'
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('E37').Value = '  +2.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01973'
$ws.Range('E38').Value = '  +2.97%  '
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '79.47'
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.454'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4555'
$ws.Range('E42').Value = '  +3.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.036'
$ws.Range('E43').Value = '  +1.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8449'
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.00'
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.804'
$ws.Range('E47').Value = '  +1.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.455'
$ws.Range('E48').Value = '  +3.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.82'
$ws.Range('E49').Value = '  +3.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4204'
$ws.Range('E50').Value = '  +3.77%  '
$ws.Range('E51').Value = '  +1.89%  '
